# Append the new match row (row 31) at the bottom of the results table,
# mirroring the existing row layout/formatting (row 30 is the template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 30
$r = 31

# Clone formatting (bold/border index on col A, datetime numFmt on col E, etc.)
# from the row above so the new row's styles stay identical to the existing ones.
$ws.Range("A$srcRow`:V$srcRow").Copy()
$ws.Range("A$r`:V$r").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($r, 1).Value = 30
$ws.Cells.Item($r, 2).Value = "india"
$ws.Cells.Item($r, 3).Value = "isl"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45233.64583333334
$ws.Cells.Item($r, 6).Value = "Odisha FC"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = "North East Utd"
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 1.41
$ws.Cells.Item($r, 11).Value = "31/10/2023 20:42"
$ws.Cells.Item($r, 12).Value = 2.14
$ws.Cells.Item($r, 13).Value = "03/11/2023 15:29"
$ws.Cells.Item($r, 14).Value = 4.92
$ws.Cells.Item($r, 15).Value = "31/10/2023 20:42"
$ws.Cells.Item($r, 16).Value = 3.69
$ws.Cells.Item($r, 17).Value = "03/11/2023 15:29"
$ws.Cells.Item($r, 18).Value = 6.78
$ws.Cells.Item($r, 19).Value = "31/10/2023 20:42"
$ws.Cells.Item($r, 20).Value = 3.26
$ws.Cells.Item($r, 21).Value = "03/11/2023 15:29"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/india/isl/odisha-fc-north-east-united/UudcwW03/"

# Make sure column E keeps the expected datetime display format.
$ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
